# Updates cryptos list prices/volumes (and three pairs of swapped rows)
# Note: Price values in column D are stored as text (e.g. "313.05"), so we
# prefix them with a leading apostrophe to force Excel to keep them as text
# instead of auto-converting to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'27.774.40"
$ws.Range("E2").Value = "  -0.56%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'1.905.42"
$ws.Range("E3").Value = "  +0.08%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.46%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'313.10"
$ws.Range("E5").Value = "  -1.16%  "

# Row 6 - USDC
$ws.Range("D6").Value = "'1.003"

# Row 7 - XRP
$ws.Range("D7").Value = "'0.4976"
$ws.Range("E7").Value = "  +3.08%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.3778"
$ws.Range("E8").Value = "  -0.41%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.07254"
$ws.Range("E9").Value = "  -1.56%  "

# Row 10 - Solana
$ws.Range("D10").Value = "'21.12"
$ws.Range("E10").Value = "  +1.70%  "

# Row 11 - Polygon
$ws.Range("D11").Value = "'0.9014"
$ws.Range("E11").Value = "  -3.29%  "

# Row 12 - TRON
$ws.Range("D12").Value = "'0.07636"
$ws.Range("E12").Value = "  -1.35%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "'1.915.53"
$ws.Range("E13").Value = "  -1.26%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'5.458"
$ws.Range("E14").Value = "  -0.48%  "

# Row 15 - Litecoin
$ws.Range("D15").Value = "'91.82"

# Row 16 - BinanceUSD
$ws.Range("E16").Value = "  -0.53%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "'0.000008709"
$ws.Range("E17").Value = "  -1.75%  "

# Row 18 - Dai
$ws.Range("D18").Value = "'1.002"
$ws.Range("E18").Value = "  -0.30%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "'27.796.80"
$ws.Range("E19").Value = "  -0.70%  "

# Row 20 - Avalanche
$ws.Range("D20").Value = "'14.54"
$ws.Range("E20").Value = "  -0.81%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +0.22%  "

# Row 22 - WrappedliquidstakedEther2.0
$ws.Range("D22").Value = "'2.161.51"
$ws.Range("E22").Value = "  -0.19%  "

# Row 23 - Cosmos
$ws.Range("E23").Value = "  -0.79%  "

# Row 24 - Chainlink
$ws.Range("D24").Value = "'6.574"
$ws.Range("E24").Value = "  -0.87%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'152.88"
$ws.Range("E25").Value = "  -2.03%  "

# Row 26 - Toncoin
$ws.Range("D26").Value = "'1.849"
$ws.Range("E26").Value = "  -3.47%  "

# Row 27 - LidoDAOToken
$ws.Range("D27").Value = "'2.212"
$ws.Range("E27").Value = "  +3.88%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'18.32"
$ws.Range("E28").Value = "  -0.82%  "

# Row 29 - BitcoinCash
$ws.Range("D29").Value = "'115.03"
$ws.Range("E29").Value = "  -1.90%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").Value = "'4.871"
$ws.Range("E30").Value = "  -1.92%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "'0.08946"
$ws.Range("E31").Value = "  +0.08%  "

# Row 32 - HuobiToken
$ws.Range("D32").Value = "'3.197"
$ws.Range("E32").Value = "  -2.02%  "

# Row 33 - ImmutableX
$ws.Range("D33").Value = "'0.7844"
$ws.Range("E33").Value = "  +2.39%  "

# Row 34 - was ARBITRUM, now Filecoin
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'4.792"
$ws.Range("E34").Value = "  +2.70%  "

# Row 35 - was Filecoin, now ARBITRUM
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.230"
$ws.Range("E35").Value = "  -1.88%  "

# Row 36 - RenderToken
$ws.Range("D36").Value = "'2.638"
$ws.Range("E36").Value = "  +3.68%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "'0.02074"
$ws.Range("E37").Value = "  +0.96%  "

# Row 38 - MXToken
$ws.Range("E38").Value = "  +1.93%  "

# Row 39 - TrustWalletToken
$ws.Range("E39").Value = "  -0.92%  "

# Row 40 - was Hedera, now TheSandbox
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.5508"
$ws.Range("E40").Value = "  +0.51%  "

# Row 41 - was TheSandbox, now Hedera
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").Value = "'0.05290"
$ws.Range("E41").Value = "  +0.28%  "

# Row 42 - FraxShare
$ws.Range("D42").Value = "'6.751"
$ws.Range("E42").Value = "  -2.78%  "

# Row 43 - Quant
$ws.Range("D43").Value = "'114.43"
$ws.Range("E43").Value = "  +4.11%  "

# Row 44 - Aptos
$ws.Range("D44").Value = "'8.461"
$ws.Range("E44").Value = "  -0.13%  "

# Row 45 - Algorand
$ws.Range("D45").Value = "'0.1510"
$ws.Range("E45").Value = "  -1.00%  "

# Row 46 - was EnergySwap, now Decentraland
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.4779"
$ws.Range("E46").Value = "  -0.53%  "

# Row 47 - was Decentraland, now EnergySwap
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'10.44"
$ws.Range("E47").Value = "  -2.19%  "

# Row 48 - PaxDollar
$ws.Range("D48").Value = "'1.003"
$ws.Range("E48").Value = "  -0.29%  "

# Row 49 - NEARProtocol
$ws.Range("D49").Value = "'1.631"
$ws.Range("E49").Value = "  -0.92%  "

# Row 50 - Aave
$ws.Range("D50").Value = "'67.08"
$ws.Range("E50").Value = "  -1.12%  "

# Row 51 - Cronos
$ws.Range("D51").Value = "'0.06019"
$ws.Range("E51").Value = "  -1.07%  "
